$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update "last updated" timestamp string
$ws.Range("A1").Value = "Datos actualizados a 28 de Abril de 2020 a las 04:52"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1010507
$ws.Range("C4").Value = 151
$ws.Range("E4").Value = 814714
$ws.Range("G4").Value = 6
$ws.Range("H4").Value = 56803

# Row 14 - Brasil
$ws.Range("B14").Value = 67446
$ws.Range("C14").Value = 945
$ws.Range("E14").Value = 31701
$ws.Range("G14").Value = 60
$ws.Range("H14").Value = 4603

# Row 30 - Japon
$ws.Range("B30").Value = 14325
$ws.Range("C30").Value = 172
$ws.Range("E30").Value = 12032
$ws.Range("G30").Value = 9
$ws.Range("H30").Value = 394

# Row 129 - Paraguay
$ws.Range("B129").Value = 230
$ws.Range("C129").Value = 2
$ws.Range("D129").Value = 95
